$wb = $excel.ActiveWorkbook

# Remove the "Texas Notes" sheet entirely - its learning-rate research notes
# are no longer part of this workbook.
$excel.DisplayAlerts = $false
$notes = $wb.Worksheets.Item("Texas Notes")
$notes.Delete()
$excel.DisplayAlerts = $true

# The PDiCECpDoC sheet used to pull its "average learning rate" from a
# formula on the Texas Notes sheet (='Texas Notes'!B10, which evaluated to
# 0.1525). Now that sheet is gone, so replace it with the plain numeric
# value that was used before that average was introduced (0.13).
$ws = $wb.Worksheets.Item("PDiCECpDoC")
$ws.Range("B2").Value = 0.13

# Leave the cursor on B2 of PDiCECpDoC ...
$ws.Activate()
$ws.Range("B2").Select()

# ... but land back on the About sheet as the one shown when the workbook
# is reopened.
$about = $wb.Worksheets.Item("About")
$about.Activate()
